# Separated enums into PdfEnums.cs file.
# Adds newly-tracked PDFium functions to the tracker sheet and marks a few
# previously-untracked rows in column P ("done" marker, matching the
# existing gray/bold/centered style already used elsewhere in that column).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Mark column P ("x") for rows 7, 9 and 10, copying the formatting ---
# --- already used by the rest of that column (e.g. P1) so the new    ---
# --- cells share the same bold/gray/centered style.                  ---
$ws.Range("P1").Copy()
$ws.Range("P7").PasteSpecial(-4122)
$ws.Range("P9").PasteSpecial(-4122)
$ws.Range("P10").PasteSpecial(-4122)

$ws.Range("P7").Value = "x"
$ws.Range("P9").Value = "x"
$ws.Range("P10").Value = "x"

# --- Newly tracked functions, listed in column O rows 18-26 ---
$ws.Range("O18").Value = "Page_HasTransparency"
$ws.Range("O19").Value = "PageObj_GetIsActive"
$ws.Range("O20").Value = "PageObj_SetIsActive"
$ws.Range("O21").Value = "PageObj_Transform"
$ws.Range("O22").Value = "PageObjTransformF"
$ws.Range("O23").Value = "PageObj_GetMatrix"
$ws.Range("O24").Value = "PageObj_SetMatrix"
$ws.Range("O25").Value = "Page_TransformAnnots"
$ws.Range("O26").Value = "PageObj_NewImageObj"

# --- Move the active selection, matching where the author was working ---
$ws.Range("O27").Select()
